# feat: update foreach/endrow/endloop with new behaviour
#
# Adds a new "#! END_ROW true" marker column (E) next to the existing
# "FOR_EACH" / "END_LOOP" template rows, mirroring the existing
# "#! END_ROW" marker on the row below it, and moves the active
# selection to H12 (matching the saved state of the authored workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column E values (row 1 gets the new "end row with loop-continuation"
# marker, row 2 reuses the existing plain "#! END_ROW" marker string).
$ws.Range("E1").Value = "#! END_ROW true"
$ws.Range("E2").Value = "#! END_ROW"

# Update the selected cell/range shown when the workbook is reopened.
$ws.Range("H12").Select()
